$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asesorias")

$ws.Range("A2").Value = "david soto"
$ws.Range("B2").Value = "Maryem Ruíz"
$ws.Range("C2").Value = "Consulta general"
$ws.Range("D2").Value = "31-10-2023"
$ws.Range("E2").Value = "00:20 - 00:40"
